$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1464.579
$ws.Range("I4").Value = 1071.3334
$ws.Range("J4").Value = 1818.5
$ws.Range("K4").Value = 1071.3334
$ws.Range("L4").Value = 1818.5
$ws.Range("M4").Value = -957.3334
$ws.Range("N4").Value = -2046.5

$ws.Range("H76").Value = 14223
$ws.Range("I76").Value = 15631
$ws.Range("J76").Value = 9999
$ws.Range("K76").Value = 15631
$ws.Range("L76").Value = 9999
$ws.Range("M76").Value = -15316
$ws.Range("N76").Value = -10629

$ws.Range("H79").Value = 14223
$ws.Range("I79").Value = 15631
$ws.Range("J79").Value = 9999
$ws.Range("K79").Value = 15631
$ws.Range("L79").Value = 9999
$ws.Range("M79").Value = -14539
$ws.Range("N79").Value = -12183

$ws.Range("H103").Value = 41668548
$ws.Range("I103").Value = 1498.5
$ws.Range("J103").Value = 50001956
$ws.Range("K103").Value = 4495.5
$ws.Range("L103").Value = 150005868
$ws.Range("M103").Value = -3909.5
$ws.Range("N103").Value = -150007040

$ws.Range("H137").Value = 41671868
$ws.Range("I137").Value = 125004490
$ws.Range("J137").Value = 5560.875
$ws.Range("K137").Value = 375013470
$ws.Range("L137").Value = 16682.625
$ws.Range("M137").Value = -375010920
$ws.Range("N137").Value = -21782.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1294.5454
$ws.Range("I2").Value = 1154.1333
$ws.Range("J2").Value = 1595.4286
$ws.Range("K2").Value = 1154.1333
$ws.Range("L2").Value = 1595.4286
$ws.Range("M2").Value = -1041.1333
$ws.Range("N2").Value = -1821.4286

$ws.Range("H97").Value = 1135.8334
$ws.Range("I97").Value = 770.4737
$ws.Range("J97").Value = 2524.2
$ws.Range("K97").Value = 770.4737
$ws.Range("L97").Value = 2524.2
$ws.Range("M97").Value = -274.4737
$ws.Range("N97").Value = -3516.2

$ws.Range("H116").Value = 1294.5454
$ws.Range("I116").Value = 1154.1333
$ws.Range("J116").Value = 1595.4286
$ws.Range("K116").Value = 1154.1333
$ws.Range("L116").Value = 1595.4286
$ws.Range("M116").Value = 1139.8667
$ws.Range("N116").Value = -6183.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1294.5454
$ws.Range("I3").Value = 1154.1333
$ws.Range("J3").Value = 1595.4286
$ws.Range("K3").Value = 1154.1333
$ws.Range("L3").Value = 1595.4286
$ws.Range("M3").Value = -1040.1333
$ws.Range("N3").Value = -1823.4286

$ws.Range("H20").Value = 6464.7417
$ws.Range("I20").Value = 9216.467000000001
$ws.Range("J20").Value = 3885
$ws.Range("K20").Value = 9216.467000000001
$ws.Range("L20").Value = 3885
$ws.Range("M20").Value = -8969.467000000001
$ws.Range("N20").Value = -4379

$ws.Range("H94").Value = 2336.5264
$ws.Range("I94").Value = 2926.0833
$ws.Range("J94").Value = 1325.8572
$ws.Range("K94").Value = 2926.0833
$ws.Range("L94").Value = 1325.8572
$ws.Range("M94").Value = -2475.0833
$ws.Range("N94").Value = -2227.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 442.7143
$ws.Range("I7").Value = 609.6667
$ws.Range("J7").Value = 317.5
$ws.Range("K7").Value = 609.6667
$ws.Range("L7").Value = 317.5
$ws.Range("M7").Value = -496.6667
$ws.Range("N7").Value = -543.5

$ws.Range("H31").Value = 23812358
$ws.Range("I31").Value = 27780114
$ws.Range("J31").Value = 5824.1665
$ws.Range("K31").Value = 27780114
$ws.Range("L31").Value = 5824.1665
$ws.Range("M31").Value = -27779819
$ws.Range("N31").Value = -6414.1665

$ws.Range("H34").Value = 23812358
$ws.Range("I34").Value = 27780114
$ws.Range("J34").Value = 5824.1665
$ws.Range("K34").Value = 27780114
$ws.Range("L34").Value = 5824.1665
$ws.Range("M34").Value = -27779912
$ws.Range("N34").Value = -6228.1665

$ws.Range("H103").Value = 35786.547
$ws.Range("I103").Value = 14335.2
$ws.Range("J103").Value = 53662.668
$ws.Range("K103").Value = 14335.2
$ws.Range("L103").Value = 53662.668
$ws.Range("M103").Value = -13163.2
$ws.Range("N103").Value = -56006.668

$ws.Range("H132").Value = 1468.45
$ws.Range("I132").Value = 1362.3529
$ws.Range("J132").Value = 2069.6667
$ws.Range("K132").Value = 4087.0587
$ws.Range("L132").Value = 6209.000100000001
$ws.Range("M132").Value = -1557.0587
$ws.Range("N132").Value = -11269.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 3490.5
$ws.Range("I26").Value = 67.833336
$ws.Range("J26").Value = 8624.5
$ws.Range("K26").Value = 203.500008
$ws.Range("L26").Value = 25873.5
$ws.Range("M26").Value = 84.49999199999999
$ws.Range("N26").Value = -26449.5

$ws.Range("H60").Value = 7532.7
$ws.Range("I60").Value = 279
$ws.Range("J60").Value = 14786.4
$ws.Range("K60").Value = 837
$ws.Range("L60").Value = 44359.2
$ws.Range("M60").Value = -586
$ws.Range("N60").Value = -44861.2

$ws.Range("H69").Value = 4211.1055
$ws.Range("I69").Value = 1161.2041
$ws.Range("J69").Value = 22891.75
$ws.Range("K69").Value = 3483.6123
$ws.Range("L69").Value = 68675.25
$ws.Range("M69").Value = -2672.6123
$ws.Range("N69").Value = -70297.25

$ws.Range("H72").Value = 4211.1055
$ws.Range("I72").Value = 1161.2041
$ws.Range("J72").Value = 22891.75
$ws.Range("K72").Value = 10450.8369
$ws.Range("L72").Value = 206025.75
$ws.Range("M72").Value = -6394.836899999998
$ws.Range("N72").Value = -214137.75

$ws.Range("H92").Value = 220.5
$ws.Range("I92").Value = 250
$ws.Range("J92").Value = 214.6
$ws.Range("K92").Value = 750
$ws.Range("L92").Value = 643.8
$ws.Range("M92").Value = 498
$ws.Range("N92").Value = -3139.8

$ws.Range("H97").Value = 313.42856
$ws.Range("I97").Value = 350
$ws.Range("J97").Value = 307.33334
$ws.Range("K97").Value = 1050
$ws.Range("L97").Value = 922.0000200000001
$ws.Range("M97").Value = -554
$ws.Range("N97").Value = -1914.00002

$ws.Range("H113").Value = 1969.3125
$ws.Range("I113").Value = 1730.125
$ws.Range("J113").Value = 2208.5
$ws.Range("K113").Value = 5190.375
$ws.Range("L113").Value = 6625.5
$ws.Range("M113").Value = -3020.375
$ws.Range("N113").Value = -10965.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 118.181816
$ws.Range("I2").Value = 142
$ws.Range("J2").Value = 76.5
$ws.Range("K2").Value = 142
$ws.Range("L2").Value = 76.5
$ws.Range("M2").Value = -29
$ws.Range("N2").Value = -302.5

$ws.Range("H15").Value = 35999.5
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 35999.5
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 35999.5
$ws.Range("N15").Value = -36575.5

$ws.Range("H21").Value = 40216400
$ws.Range("I21").Value = 50263500
$ws.Range("J21").Value = 27997
$ws.Range("K21").Value = 50263500
$ws.Range("L21").Value = 27997
$ws.Range("M21").Value = -50263327
$ws.Range("N21").Value = -28343

$ws.Range("H30").Value = 40216400
$ws.Range("I30").Value = 50263500
$ws.Range("J30").Value = 27997
$ws.Range("K30").Value = 50263500
$ws.Range("L30").Value = 27997
$ws.Range("M30").Value = -50263395
$ws.Range("N30").Value = -28207

$ws.Range("H81").Value = 35999.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 35999.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 35999.5
$ws.Range("N81").Value = -37995.5

$ws.Range("H84").Value = 35999.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 35999.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 107998.5
$ws.Range("N84").Value = -117982.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H46").Value = 1194.9474
$ws.Range("I46").Value = 989.25
$ws.Range("J46").Value = 1344.5454
$ws.Range("K46").Value = 989.25
$ws.Range("L46").Value = 1344.5454
$ws.Range("M46").Value = -801.25
$ws.Range("N46").Value = -1720.5454

$ws.Range("H80").Value = 46394.8
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 46394.8
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 46394.8
$ws.Range("N80").Value = -48640.8

$ws.Range("H82").Value = 2642.3125
$ws.Range("I82").Value = 702.5238000000001
$ws.Range("J82").Value = 6345.5454
$ws.Range("K82").Value = 702.5238000000001
$ws.Range("L82").Value = 6345.5454
$ws.Range("M82").Value = -341.5238000000001
$ws.Range("N82").Value = -7067.5454

$ws.Range("H83").Value = 46394.8
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 46394.8
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 139184.4
$ws.Range("N83").Value = -150416.4

$ws.Range("H85").Value = 2642.3125
$ws.Range("I85").Value = 702.5238000000001
$ws.Range("J85").Value = 6345.5454
$ws.Range("K85").Value = 702.5238000000001
$ws.Range("L85").Value = 6345.5454
$ws.Range("M85").Value = 545.4761999999999
$ws.Range("N85").Value = -8841.545399999999

$ws.Range("H132").Value = 3823.182
$ws.Range("I132").Value = 2729.25
$ws.Range("J132").Value = 6740.3335
$ws.Range("K132").Value = 8187.75
$ws.Range("L132").Value = 20221.0005
$ws.Range("M132").Value = -5657.75
$ws.Range("N132").Value = -25281.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1640.1428
$ws.Range("I100").Value = 851.0909
$ws.Range("J100").Value = 4533.3335
$ws.Range("K100").Value = 1702.1818
$ws.Range("L100").Value = 9066.666999999999
$ws.Range("M100").Value = -1161.1818
$ws.Range("N100").Value = -10148.667

$ws.Range("H122").Value = 1766.8529
$ws.Range("I122").Value = 1764.1111
$ws.Range("J122").Value = 1777.4286
$ws.Range("K122").Value = 5292.3333
$ws.Range("L122").Value = 5332.2858
$ws.Range("M122").Value = -2842.3333
$ws.Range("N122").Value = -10232.2858

$ws.Range("H132").Value = 771631.4399999999
$ws.Range("I132").Value = 2721.5
$ws.Range("J132").Value = 3334664.8
$ws.Range("K132").Value = 8164.5
$ws.Range("L132").Value = 10003994.4
$ws.Range("M132").Value = -5634.5
$ws.Range("N132").Value = -10009054.4
